$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card15")

# Clear the "Serviced by" value for row 2 (was "م.صيام")
$ws.Range("O2").Value = ""

# Set the "Serviced by" value for row 3 to the literal text "nan"
$ws.Range("O3").Value = "nan"
